# First pass at adding POIs (Points of Interest) to the town map.
# Adds two new columns (Q, R) for state-level historical-marker metadata:
#   StateHistoricalMarkerFileName / StateHistoricalMarkerOneDriveLink
# and populates them for the New Hampshire row (row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) -- written first so the shared-string table
# picks up the header text before the New Hampshire data values.
$ws.Range("Q1").Value = "StateHistoricalMarkerFileName"
$ws.Range("R1").Value = "StateHistoricalMarkerOneDriveLink"

# New data cells for the New Hampshire row (row 9).
$ws.Range("Q9").Value = "NH Historical Highway Markers.2.xlsx"
$ws.Range("R9").Value = "1drv.ms/x/s!An0k-SnslkINzx1beYrxBM-Rlm6j?e=mS31Zx"

# Size the two new columns to fit their contents (matches the widened
# "N" column family already used for file-name / link columns).
$ws.Columns.Item(17).ColumnWidth = 30.166666666666668
$ws.Columns.Item(18).ColumnWidth = 28.756510416666668

# Move the selection to where the author left off editing.
$ws.Range("O8").Select()
